# "Allow null value on importer" — add three more sample rows to "Sheet 1"
# (geom_code=BC / indicator_shortcode=IND2 / admin_level=2 / Parents="B,Top")
# whose "value" column is intentionally left blank, and relax the COUNTIFS
# formulas on the COUNT sheets so a blank value isn't counted as a record.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("Sheet 1")

$rows = @(20, 21, 22)
foreach ($r in $rows) {
    $data.Range("A$r").Value = "BC"
    $data.Range("B$r").Value = "IND2"
    # Column C (value) intentionally left blank -- this is the null value case.
    # Touch a no-op font property so the (empty) cell is still materialised.
    $data.Range("C$r").Font.Underline = $False
    $data.Range("D$r").Value = 40179
    $data.Range("E$r").Value = 0
    $data.Range("F$r").Value = 40179
    $data.Range("G$r").Value = 2
    $data.Range("H$r").Value = "B,Top"
}

# COUNT: exclude rows with a blank value from the per-indicator counts.
$count = $wb.Worksheets.Item("COUNT")
for ($r = 2; $r -le 7; $r++) {
    $count.Range("D$r").Formula = "=COUNTIFS('Sheet 1'!D:D,B$r,'Sheet 1'!A:A,A$r,'Sheet 1'!B:B,C$r,'Sheet 1'!C:C,""<>""&"""")"
}

# COUNT (Upper Level): same relaxation for the rolled-up counts.
$countUpper = $wb.Worksheets.Item("COUNT (Upper Level)")
for ($r = 2; $r -le 4; $r++) {
    $countUpper.Range("D$r").Formula = "=COUNTIFS('Sheet 1'!D:D,B$r,'Sheet 1'!B:B,C$r,'Sheet 1'!H:H,""*""&A$r&""*"",'Sheet 1'!C:C,""<>""&"""")"
}
